$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 201; this shifts the existing rows 201-253
# down to 202-254 (and extends the sheet dimension to A1:R254).
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new record.
$ws.Cells.Item(201, 1).Value2 = 5
$ws.Cells.Item(201, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(201, 3).Value2 = "Maule"
$ws.Cells.Item(201, 4).Value2 = 44889
$ws.Cells.Item(201, 5).Value2 = 7
$ws.Cells.Item(201, 6).Value2 = 100112021
$ws.Cells.Item(201, 7).Value2 = "Ají"
$ws.Cells.Item(201, 8).Value2 = "Americana (o)"
$ws.Cells.Item(201, 9).Value2 = "Primera"
$ws.Cells.Item(201, 10).Value2 = 150
$ws.Cells.Item(201, 11).Value2 = 15000
$ws.Cells.Item(201, 12).Value2 = 15000
$ws.Cells.Item(201, 13).Value2 = 15000
$ws.Cells.Item(201, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(201, 15).Value2 = "Región del Maule"
$ws.Cells.Item(201, 16).Value2 = 1000
$ws.Cells.Item(201, 17).Value2 = 15
$ws.Cells.Item(201, 18).Value2 = "Hortaliza"
